# Insert a new data row at row 104 (pushes existing rows 104:151 down to 105:152)
# and populate the new row with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(104).Insert()

$ws.Cells.Item(104, 1).Value = 7
$ws.Cells.Item(104, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(104, 3).Value = "Ñuble"
$ws.Cells.Item(104, 4).Value = 44523
$ws.Cells.Item(104, 5).Value = 16
$ws.Cells.Item(104, 6).Value = 100112003
$ws.Cells.Item(104, 7).Value = "Ajo"
$ws.Cells.Item(104, 8).Value = "Chino"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 100
$ws.Cells.Item(104, 11).Value = 17000
$ws.Cells.Item(104, 12).Value = 18000
$ws.Cells.Item(104, 13).Value = 17500
$ws.Cells.Item(104, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(104, 15).Value = "China"
$ws.Cells.Item(104, 16).Value = 1750
$ws.Cells.Item(104, 17).Value = 10
$ws.Cells.Item(104, 18).Value = "Hortaliza"
